$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.4
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.25
$ws.Range("X2").Value = 11
$ws.Range("AK2").Value = 34
$ws.Range("AL2").Value = 26
$ws.Range("AR2").Value = 67
$ws.Range("AY2").Value = 26
$ws.Range("AZ2").Value = 51

$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("X3").Value = 6
$ws.Range("Z3").Value = 9
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 501
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 41
$ws.Range("AJ3").Value = 23
$ws.Range("AW3").Value = 8.5
$ws.Range("BA3").Value = 201
$ws.Range("BB3").Value = 401

$ws.Range("S4").Value = 1.24
$ws.Range("T4").Value = 3.91

$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 3.45
$ws.Range("J5").Value = 2.57
$ws.Range("K5").Value = 2.07
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 10
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.55
$ws.Range("V5").Value = 1.98
$ws.Range("W5").Value = 7.8
$ws.Range("X5").Value = 10
$ws.Range("Z5").Value = 18.5
$ws.Range("AA5").Value = 15.5
$ws.Range("AB5").Value = 25
$ws.Range("AC5").Value = 10.25
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 60
$ws.Range("AG5").Value = 450
$ws.Range("AH5").Value = 10.25
$ws.Range("AI5").Value = 18.5
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 50
$ws.Range("AL5").Value = 32
$ws.Range("AM5").Value = 37
$ws.Range("AN5").Value = 3.9
$ws.Range("AO5").Value = 10
$ws.Range("AQ5").Value = 37
$ws.Range("AR5").Value = 65
$ws.Range("AS5").Value = 250
$ws.Range("AT5").Value = 2.52
$ws.Range("AU5").Value = 7
$ws.Range("AW5").Value = 5.3
$ws.Range("AX5").Value = 19.5

$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
